$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    "hören",
    "brauchen",
    "sprechen",
    "schämen",
    "fischen",
    "hauen",
    "leisten",
    "platzen",
    "wirken",
    "sperren",
    "werden",
    "reizen",
    "altern",
    "wärmen",
    "schließen",
    "fließen",
    "bitten",
    "meinen",
    "kümmern",
    "garen",
    "enden",
    "betteln",
    "jagen",
    "wehtun",
    "reisen",
    "lassen",
    "lehnen",
    "loben",
    "sitzen",
    "knien",
    "schrecken",
    "schneiden",
    "ärgern",
    "feiern",
    "lügen",
    "duschen",
    "planen",
    "freuen",
    "rauchen",
    "fahren",
    "sichern",
    "träumen",
    "trauen",
    "ändern",
    "hacken",
    "schalten",
    "treten",
    "mögen",
    "zögern",
    "wenden",
    "dienen",
    "warnen",
    "segnen",
    "pressen",
    "spielen",
    "pfeifen",
    "bergen",
    "münzen",
    "krachen",
    "werfen",
    "schmecken",
    "fangen",
    "hupen",
    "biegen",
    "kehren",
    "mauern",
    "schulden",
    "treiben",
    "morden",
    "stoppen",
    "schnellen",
    "gelten",
    "töten",
    "sorgen",
    "brauen",
    "wehren",
    "klagen",
    "rasen",
    "grüßen",
    "zeigen",
    "pflegen",
    "regnen",
    "malen",
    "schreiten",
    "haben",
    "bauen",
    "folgen",
    "bieten",
    "liefern",
    "formen",
    "zeichnen",
    "zielen",
    "graben",
    "jubeln",
    "quälen",
    "spinnen",
    "nutzen",
    "spüren",
    "albern",
    "sprengen",
    "herrschen",
    "liegen",
    "ehren",
    "merken",
    "läuten",
    "schlucken",
    "klettern",
    "knarren",
    "nerven",
    "schauen",
    "tropfen",
    "sterben",
    "greifen",
    "reiten",
    "lockern",
    "scheitern",
    "irren",
    "heulen",
    "tauschen",
    "erben",
    "trennen",
    "flehen",
    "heilen",
    "schwören",
    "zünden",
    "fassen",
    "siegen",
    "hoffen",
    "pflanzen",
    "schenken",
    "streichen",
    "achten",
    "gründen",
    "geben",
    "weinen",
    "machen",
    "runden",
    "sinken",
    "stellen",
    "lesen",
    "scheinen",
    "küssen",
    "drehen",
    "äußern",
    "ruhen",
    "schwächen",
    "bellen",
    "posten",
    "rufen",
    "heben",
    "saufen",
    "kichern",
    "backen",
    "strahlen",
    "führen",
    "fallen",
    "stammen",
    "schicken",
    "schlafen",
    "wetten",
    "wagen",
    "bluten",
    "filmen",
    "flüchten",
    "tanzen",
    "lenken",
    "scheiden",
    "rauben",
    "wundern",
    "boxen",
    "schlagen",
    "klingen",
    "kosten",
    "seufzen",
    "stehlen",
    "wüten",
    "passen",
    "wachsen",
    "schätzen",
    "fällen",
    "dringen",
    "decken",
    "wechseln",
    "räumen",
    "schwingen",
    "arten",
    "kämpfen",
    "stecken",
    "weichen",
    "helfen",
    "suchen",
    "tollen"
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value2 = $words[$i]
}

